$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("index_returns")

# Drop the now-unused F:G columns (old QTD/YTD overflow) by deleting them,
# which collapses the used range down to columns A:E.
$ws.Range("F1:G36").Delete()

# The "1 WEEK/MTD/QTD/YTD" return figures are stored as literal text
# (e.g. "3.31%"), not numbers -- pre-format the destination cells as Text
# so assigning the percent-looking strings does not get auto-converted
# into numeric percentage values.
$ws.Range("B4:E24").NumberFormat = "@"
$ws.Range("B26:E33").NumberFormat = "@"
$ws.Range("B35:E40").NumberFormat = "@"

# Row 1
$ws.Range("A1").Value = "INDEX RETURNS"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"
$ws.Range("D1").Value = "Unnamed: 3"
$ws.Range("E1").Value = "Unnamed: 4"

# Row 2
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "1 WEEK"
$ws.Range("C2").Value = "MTD"
$ws.Range("D2").Value = "QTD"
$ws.Range("E2").Value = "YTD"

# Row 3
$ws.Range("A3").Value = "EQUITIES"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""

# Row 4
$ws.Range("A4").Value = "S&P 500"
$ws.Range("B4").Value = "3.31%"
$ws.Range("C4").Value = "3.71%"
$ws.Range("D4").Value = "13.27%"
$ws.Range("E4").Value = "-9.34%"

# Row 5
$ws.Range("A5").Value = "DJ Industrial Average"
$ws.Range("B5").Value = "2.99%"
$ws.Range("C5").Value = "2.87%"
$ws.Range("D5").Value = "9.89%"
$ws.Range("E5").Value = "-5.98%"

# Row 6
$ws.Range("A6").Value = "Russell 2000"
$ws.Range("B6").Value = "4.97%"
$ws.Range("C6").Value = "7.03%"
$ws.Range("D6").Value = "18.21%"
$ws.Range("E6").Value = "-9.48%"

# Row 7
$ws.Range("A7").Value = "Russell Midcap"
$ws.Range("B7").Value = "4.19%"
$ws.Range("C7").Value = "4.91%"
$ws.Range("D7").Value = "15.26%"
$ws.Range("E7").Value = "-9.60%"

# Row 8
$ws.Range("A8").Value = "STOXX Europe 50 (€)"
$ws.Range("B8").Value = "1.42%"
$ws.Range("C8").Value = "1.92%"
$ws.Range("D8").Value = "9.53%"
$ws.Range("E8").Value = "-9.58%"

# Row 9
$ws.Range("A9").Value = "STOXX Europe 600 (€)†"
$ws.Range("B9").Value = "1.29%"
$ws.Range("C9").Value = "0.74%"
$ws.Range("D9").Value = "8.54%"
$ws.Range("E9").Value = "-7.71%"

# Row 10
$ws.Range("A10").Value = "MSCI EAFE Small Cap"
$ws.Range("B10").Value = "2.99%"
$ws.Range("C10").Value = "2.47%"
$ws.Range("D10").Value = "9.25%"
$ws.Range("E10").Value = "-17.47%"

# Row 11
$ws.Range("A11").Value = "FTSE 100 (£)"
$ws.Range("B11").Value = "1.18%"
$ws.Range("C11").Value = "1.54%"
$ws.Range("D11").Value = "5.27%"
$ws.Range("E11").Value = "4.25%"

# Row 12
$ws.Range("A12").Value = "DAX (€)"
$ws.Range("B12").Value = "1.63%"
$ws.Range("C12").Value = "2.31%"
$ws.Range("D12").Value = "7.92%"
$ws.Range("E12").Value = "-13.15%"

# Row 13
$ws.Range("A13").Value = "FTSE MIB (€)"
$ws.Range("B13").Value = "1.70%"
$ws.Range("C13").Value = "2.52%"
$ws.Range("D13").Value = "8.35%"
$ws.Range("E13").Value = "-12.89%"

# Row 14
$ws.Range("A14").Value = "CAC 40 (€)†"
$ws.Range("B14").Value = "1.26%"
$ws.Range("C14").Value = "1.63%"
$ws.Range("D14").Value = "10.76%"
$ws.Range("E14").Value = "-6.45%"

# Row 15
$ws.Range("A15").Value = "SWISS MKT (CHF)"
$ws.Range("B15").Value = "0.05%"
$ws.Range("C15").Value = "-0.16%"
$ws.Range("D15").Value = "3.60%"
$ws.Range("E15").Value = "-11.22%"

# Row 16
$ws.Range("A16").Value = "TOPIX (¥)"
$ws.Range("B16").Value = "1.34%"
$ws.Range("C16").Value = "1.69%"
$ws.Range("D16").Value = "5.48%"
$ws.Range("E16").Value = "0.44%"

# Row 17
$ws.Range("A17").Value = "Hang Seng (HKD)"
$ws.Range("B17").Value = "-0.11%"
$ws.Range("C17").Value = "0.13%"
$ws.Range("D17").Value = "-7.20%"
$ws.Range("E17").Value = "-11.67%"

# Row 18
$ws.Range("A18").Value = "MSCI World"
$ws.Range("B18").Value = "3.06%"
$ws.Range("C18").Value = "3.30%"
$ws.Range("D18").Value = "11.54%"
$ws.Range("E18").Value = "-11.09%"

# Row 19
$ws.Range("A19").Value = "MSCI China Free†"
$ws.Range("B19").Value = "-0.38%"
$ws.Range("C19").Value = "-0.01%"
$ws.Range("D19").Value = "-8.93%"
$ws.Range("E19").Value = "-15.12%"

# Row 20
$ws.Range("A20").Value = "MSCI EAFE"
$ws.Range("B20").Value = "2.17%"
$ws.Range("C20").Value = "1.51%"
$ws.Range("D20").Value = "6.58%"
$ws.Range("E20").Value = "-13.95%"

# Row 21
$ws.Range("A21").Value = "MSCI EM"
$ws.Range("B21").Value = "1.66%"
$ws.Range("C21").Value = "2.65%"
$ws.Range("D21").Value = "2.49%"
$ws.Range("E21").Value = "-15.42%"

# Row 22
$ws.Range("A22").Value = "MSCI Brazil (BRL)"
$ws.Range("B22").Value = "6.16%"
$ws.Range("C22").Value = "9.11%"
$ws.Range("D22").Value = "14.45%"
$ws.Range("E22").Value = "10.53%"

# Row 23
$ws.Range("A23").Value = "MSCI India (INR)"
$ws.Range("B23").Value = "1.69%"
$ws.Range("C23").Value = "3.48%"
$ws.Range("D23").Value = "13.60%"
$ws.Range("E23").Value = "2.48%"

# Row 24
$ws.Range("A24").Value = "MSCI Russia (RUB)"
$ws.Range("B24").Value = "0.00%"
$ws.Range("C24").Value = "0.00%"
$ws.Range("D24").Value = "0.00%"
$ws.Range("E24").Value = "-100.00%"

# Row 25
$ws.Range("A25").Value = "FIXED INCOME"
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = ""

# Row 26
$ws.Range("A26").Value = "Bloomberg Aggregate"
$ws.Range("B26").Value = "0.24%"
$ws.Range("C26").Value = "-0.80%"
$ws.Range("D26").Value = "1.63%"
$ws.Range("E26").Value = "-8.89%"

# Row 27
$ws.Range("A27").Value = "Bloomberg Euro Aggregate"
$ws.Range("B27").Value = "0.32%"
$ws.Range("C27").Value = "-0.51%"
$ws.Range("D27").Value = "1.05%"
$ws.Range("E27").Value = "-18.37%"

# Row 28
$ws.Range("A28").Value = "Bloomberg US High Yield"
$ws.Range("B28").Value = "0.94%"
$ws.Range("C28").Value = "1.60%"
$ws.Range("D28").Value = "7.60%"
$ws.Range("E28").Value = "-7.67%"

# Row 29
$ws.Range("A29").Value = "Bloomberg Euro High Yield (€)"
$ws.Range("B29").Value = "0.75%"
$ws.Range("C29").Value = "1.66%"
$ws.Range("D29").Value = "6.83%"
$ws.Range("E29").Value = "-8.59%"

# Row 30
$ws.Range("A30").Value = "Bloomberg Muni Aggregate"
$ws.Range("B30").Value = "-0.08%"
$ws.Range("C30").Value = "-0.19%"
$ws.Range("D30").Value = "2.44%"
$ws.Range("E30").Value = "-6.76%"

# Row 31
$ws.Range("A31").Value = "Bloomberg TIPS"
$ws.Range("B31").Value = "0.18%"
$ws.Range("C31").Value = "-1.41%"
$ws.Range("D31").Value = "3.06%"
$ws.Range("E31").Value = "-6.90%"

# Row 32
$ws.Range("A32").Value = "JPM EMBI Glbl. Divers."
$ws.Range("B32").Value = "1.51%"
$ws.Range("C32").Value = "2.43%"
$ws.Range("D32").Value = "5.39%"
$ws.Range("E32").Value = "-16.01%"

# Row 33
$ws.Range("A33").Value = "JPM GBI-EM Glbl. Divers."
$ws.Range("B33").Value = "2.34%"
$ws.Range("C33").Value = "3.09%"
$ws.Range("D33").Value = "3.40%"
$ws.Range("E33").Value = "-11.63%"

# Row 34
$ws.Range("A34").Value = "OTHER"
$ws.Range("B34").Value = ""
$ws.Range("C34").Value = ""
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = ""

# Row 35
$ws.Range("A35").Value = "DJ US Real Estate"
$ws.Range("B35").Value = "4.28%"
$ws.Range("C35").Value = "2.66%"
$ws.Range("D35").Value = "11.77%"
$ws.Range("E35").Value = "-10.61%"

# Row 36
$ws.Range("A36").Value = "FTSE EPRA/NAREIT Dvlpd. Ex-US"
$ws.Range("B36").Value = "2.89%"
$ws.Range("C36").Value = "0.83%"
$ws.Range("D36").Value = "7.21%"
$ws.Range("E36").Value = "-14.36%"

# Row 37
$ws.Range("A37").Value = "S&P GSCI"
$ws.Range("B37").Value = "4.45%"
$ws.Range("C37").Value = "-1.83%"
$ws.Range("D37").Value = "-1.87%"
$ws.Range("E37").Value = "33.26%"

# Row 38
$ws.Range("A38").Value = "Alerian MLP *"
$ws.Range("B38").Value = "4.95%"
$ws.Range("C38").Value = "2.14%"
$ws.Range("D38").Value = "14.90%"
$ws.Range("E38").Value = "26.43%"

# Row 39
$ws.Range("A39").Value = "US Dollar Index"
$ws.Range("B39").Value = "-0.91%"
$ws.Range("C39").Value = "-0.25%"
$ws.Range("D39").Value = "0.90%"
$ws.Range("E39").Value = "10.59%"

# Row 40
$ws.Range("A40").Value = "VIX"
$ws.Range("B40").Value = "-7.66%"
$ws.Range("C40").Value = "-8.44%"
$ws.Range("D40").Value = "-31.97%"
$ws.Range("E40").Value = "13.41%"

